$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("其他JS脚本")

# --- Split the combined "虎牙,斗鱼" row into two separate rows -----------------
# Insert a new row right after the existing "虎牙,斗鱼" row (row 4), then copy
# the formatting of row 4 down onto it so the new row matches the sheet's style.
$ws.Rows.Item(5).Insert()

$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight()

# Row 4 becomes the "虎牙" (Huya) entry; it now also has a PHP script.
$ws.Range("A4").Value = "虎牙"
$ws.Range("B4").Value = "HuyaAndDouyu.js"
$ws.Range("C4").Value = "huya.php"
$ws.Range("D4").Value = "虎牙平台直播间"

# Row 5 becomes the "斗鱼" (Douyu) entry, sharing the same JS script.
$ws.Range("A5").Value = "斗鱼"
$ws.Range("B5").Value = "HuyaAndDouyu.js"
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = "斗鱼平台直播间"

# --- Re-point the worksheet AutoFilter so it still ends right after the new --
# --- row, rather than auto-expanding across the rows further down. ----------
$a7 = $ws.Range("A7").Value()
$c7 = $ws.Range("C7").Value()
$d7 = $ws.Range("D7").Value()
$a8 = $ws.Range("A8").Value()
$c8 = $ws.Range("C8").Value()
$d8 = $ws.Range("D8").Value()

$ws.Range("A7:D8").ClearContents()

$ws.AutoFilterMode = $false
$ws.Range("A1:D6").AutoFilter()

$ws.Range("A7").Value = $a7
$ws.Range("C7").Value = $c7
$ws.Range("D7").Value = $d7
$ws.Range("A8").Value = $a8
$ws.Range("C8").Value = $c8
$ws.Range("D8").Value = $d8

# --- Keep the workbook-level hidden _FilterDatabase name for this sheet in --
# --- sync with the new AutoFilter range. ------------------------------------
foreach ($n in $wb.Names) {
    if ($n.RefersTo().IndexOf("其他JS脚本") -ge 0) {
        $n.RefersTo = "=其他JS脚本!`$A`$1:`$D`$6"
    }
}
